$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Activate this sheet so it becomes the saved active tab (index 4) and
# the "tabSelected" flag moves off "NewLoanInput" onto this sheet.
$ws.Activate()

# Insert a new row before the old row 8 ("clickonsubmit"/"Submit"),
# shifting the existing rows 8-14 down to 9-15.
$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = "waittopageload1"
$ws.Range("B8").Value = 2000

# Match the numeric style used by the sibling "waittopageload" row (B3)
# instead of the text style that Insert() copied down from row 7.
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on the newly inserted row, as in the target file.
$ws.Range("A8:B8").Select()
